# Applies the "Updated cryptos list" snapshot refresh (prices + 1h volume %,
# plus two ranking swaps) produced by the Sun Dec 24 16:37:03 UTC 2023
# GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '43.806.20'
$ws.Range("E2").Value = '  +0.01%  '

$ws.Range("D3").Value = "'" + '2.293.72'
$ws.Range("E3").Value = '  -0.05%  '

$ws.Range("E4").Value = '  +0.17%  '

$ws.Range("D5").Value = "'" + '114.40'
$ws.Range("E5").Value = '  +17.57%  '

$ws.Range("D6").Value = "'" + '269.65'
$ws.Range("E6").Value = '  -0.08%  '

$ws.Range("E7").Value = '  +0.30%  '

$ws.Range("E8").Value = '  +0.27%  '

$ws.Range("E9").Value = '  +1.31%  '

$ws.Range("D10").Value = "'" + '48.56'
$ws.Range("E10").Value = '  +7.06%  '

$ws.Range("D11").Value = "'" + '0.0945'
$ws.Range("E11").Value = '  +1.43%  '

$ws.Range("D12").Value = "'" + '8.93'
$ws.Range("E12").Value = '  +13.47%  '

$ws.Range("D14").Value = "'" + '15.83'
$ws.Range("E14").Value = '  -0.18%  '

$ws.Range("D15").Value = "'" + '2.635.59'
$ws.Range("E15").Value = '  +0.03%  '

$ws.Range("D16").Value = "'" + '0.858'
$ws.Range("E16").Value = '  -0.31%  '

$ws.Range("D17").Value = "'" + '2.302.10'
$ws.Range("E17").Value = '  +0.81%  '

$ws.Range("D18").Value = "'" + '43.703.51'
$ws.Range("E18").Value = '  -0.22%  '

$ws.Range("E19").Value = '  -0.61%  '

$ws.Range("D20").Value = "'" + '6.93'
$ws.Range("E20").Value = '  +11.80%  '

$ws.Range("D21").Value = "'" + '72.29'
$ws.Range("E21").Value = '  +0.12%  '

$ws.Range("E22").Value = '  -0.92%  '

$ws.Range("D23").Value = "'" + '9.86'
$ws.Range("E23").Value = '  +8.52%  '

$ws.Range("D24").Value = "'" + '232.81'
$ws.Range("E24").Value = '  -0.12%  '

$ws.Range("E25").Value = '  +8.37%  '

$ws.Range("E26").Value = '  +0.04%  '

$ws.Range("D27").Value = "'" + '11.58'
$ws.Range("E27").Value = '  +2.62%  '

$ws.Range("E28").Value = '  -0.87%  '

$ws.Range("D29").Value = "'" + '42.18'
$ws.Range("E29").Value = '  +9.68%  '

$ws.Range("E30").Value = '  -2.01%  '

$ws.Range("E31").Value = '  -1.11%  '

$ws.Range("D32").Value = "'" + '175.58'
$ws.Range("E32").Value = '  -0.38%  '

$ws.Range("D33").Value = "'" + '21.61'
$ws.Range("E33").Value = '  -1.03%  '

$ws.Range("E34").Value = '  +4.04%  '

$ws.Range("D35").Value = "'" + '5.71'
$ws.Range("E35").Value = '  +5.27%  '

$ws.Range("E36").Value = '  -0.09%  '

$ws.Range("D37").Value = "'" + '4.75'
$ws.Range("E37").Value = '  +1.09%  '

$ws.Range("D38").Value = "'" + '0.0363'
$ws.Range("E38").Value = '  +3.41%  '

$ws.Range("E39").Value = '  -2.83%  '

$ws.Range("D40").Value = "'" + '3.83'
$ws.Range("E40").Value = '  +10.48%  '

$ws.Range("D41").Value = "'" + '13.94'
$ws.Range("E41").Value = '  +14.10%  '

$ws.Range("D42").Value = "'" + '74.29'
$ws.Range("E42").Value = '  +15.31%  '

$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").Value = "'" + '0.243'
$ws.Range("E43").Value = '  +2.36%  '

$ws.Range("B44").Value = 'LidoDAOToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D44").Value = "'" + '2.40'
$ws.Range("E44").Value = '  +2.92%  '

$ws.Range("D45").Value = "'" + '6.38'
$ws.Range("E45").Value = '  +22.15%  '

$ws.Range("E46").Value = '  +0.08%  '

$ws.Range("E47").Value = '  +3.13%  '

$ws.Range("E48").Value = '  +0.90%  '

$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").Value = "'" + '102.60'
$ws.Range("E49").Value = '  +3.51%  '

$ws.Range("B50").Value = 'TrustWalletToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D50").Value = "'" + '1.26'
$ws.Range("E50").Value = '  +2.20%  '

$ws.Range("D51").Value = "'" + '0.0994'
$ws.Range("E51").Value = '  -2.97%  '
